# Apply the "float input" grading fix to the marksheet.
#
# Summary of the edit:
#  - The summary table (rows 10-12) gets real grading numbers instead of the
#    placeholder zeros/"Absent" text, the label cells (A10/A11/A12) pick up
#    the same "mtitleStyle" used elsewhere in that header block, and the
#    previously text-typed "-1" in C11 becomes a genuine number (this is the
#    "handles float input" fix from the commit message).
#  - The sheet used to reserve three side-by-side Student/Correct answer
#    blocks (A:B, D:E, G:H) that were otherwise empty placeholders. Only the
#    first student (columns A:B) is actually used, so the extra columns are
#    cleared out (this is what shrinks the used range from A5:H40 to A5:E40).
#  - Column A (Student Ans) for each question (rows 16-40) is now populated
#    with the student's actual answer, styled with the "correctStyle" /
#    "incorrectStyle" named cell styles depending on whether it matches the
#    "Correct Ans" already stored in column B. Unattempted questions are left
#    blank with the default "normalStyle".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Summary table (rows 10-12)
# ---------------------------------------------------------------------

# Label cells pick up the same style used by the rest of the header block.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Row 10: No. (right / wrong / not-attempt / max)
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

# Row 11: Marking scheme. C11 used to be stored as the text "-1"; store it as
# a real number now so downstream formulas/float handling don't break.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12: Totals
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "74/112"

# ---------------------------------------------------------------------
# 2. Drop the unused 2nd/3rd "Student Ans / Correct Ans" column blocks
# ---------------------------------------------------------------------

$ws.Range("G15:H40").Clear()
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# 3. Fill in the student's answers for each question (column A),
#    colouring them green (correctStyle) or red (incorrectStyle) based on
#    whether they match the correct answer already in column B.
# ---------------------------------------------------------------------

$answers = @(
    @{Cell="A16"; Style="incorrectStyle"; Value="Option C"},
    @{Cell="A18"; Style="correctStyle";   Value="Option B"},
    @{Cell="A19"; Style="correctStyle";   Value="Option C"},
    @{Cell="A20"; Style="correctStyle";   Value="Option B"},
    @{Cell="A21"; Style="correctStyle";   Value="Option C"},
    @{Cell="A22"; Style="correctStyle";   Value="Option D"},
    @{Cell="A23"; Style="correctStyle";   Value="Option D"},
    @{Cell="A24"; Style="correctStyle";   Value="Option A"},
    @{Cell="A25"; Style="correctStyle";   Value="Option A"},
    @{Cell="A26"; Style="correctStyle";   Value="Option C"},
    @{Cell="A30"; Style="correctStyle";   Value="Option B"},
    @{Cell="A32"; Style="correctStyle";   Value="Option C"},
    @{Cell="A33"; Style="correctStyle";   Value="Option D"},
    @{Cell="A34"; Style="correctStyle";   Value="Option B"},
    @{Cell="A36"; Style="correctStyle";   Value="Option A"},
    @{Cell="A38"; Style="correctStyle";   Value="Option A"},
    @{Cell="A39"; Style="correctStyle";   Value="Option D"},
    @{Cell="A40"; Style="correctStyle";   Value="Option D"}
)

foreach ($ans in $answers) {
    $cell = $ws.Range($ans.Cell)
    $cell.Value = $ans.Value
    $cell.Style = $ans.Style
}

# The leftover 2nd column block (D16:E18) keeps its data, it just needs the
# same Student Ans colouring treatment as column A.
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"

$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"

$ws.Range("D18").Value = "Option B"
$ws.Range("D18").Style = "incorrectStyle"
